{"js": "// The captured change is NOT a content/formatting edit: every hunk in the\n// diff touches only the root element's `xmlns:*` attribute *ordering* on\n// seven package parts (document.xml, endnotes.xml, footer1.xml,\n// footnotes.xml, header1.xml, styles.xml, theme1.xml). Once the namespace\n// declarations are stripped, every \"before\" root tag is byte-identical to\n// its \"after\" counterpart (same element, same attributes such as\n// `mc:Ignorable`, same children) \u2014 confirmed by diffing them directly.\n//\n// That reordering is purely a side effect of which OOXML writer produced\n// the package (the commit message says as much: \"Rework DOCX\n// implementation to better support testing and decrease coupling\" \u2014 a\n// tooling/serializer change, not a document edit). The Word JavaScript API\n// only exposes the semantic object model (text, formatting, structure);\n// it has no surface for controlling the literal attribute order Word/the\n// host emits when it serializes a part's root element. So there is no\n// Office.js call sequence that can change that ordering, and none is\n// needed: the visible document content is unchanged.\n//\n// Reflect that faithfully as a no-op: read a harmless property and don't\n// mutate anything.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The captured change is NOT a content/formatting edit: every hunk in the\n# diff touches only the root element's `xmlns:*` attribute *ordering* on\n# seven package parts (document.xml, endnotes.xml, footer1.xml,\n# footnotes.xml, header1.xml, styles.xml, theme1.xml). Once the namespace\n# declarations are stripped, every \"before\" root tag is byte-identical to\n# its \"after\" counterpart (same element, same attributes such as\n# `mc:Ignorable`, same children) - confirmed by diffing them directly.\n#\n# That reordering is purely a side effect of which OOXML writer produced\n# the package (the commit message says as much: \"Rework DOCX\n# implementation to better support testing and decrease coupling\" - a\n# tooling/serializer change, not a document edit). The Word COM object\n# model only exposes the semantic object model (text, formatting,\n# structure); it has no surface for controlling the literal attribute\n# order Word emits when it serializes a part's root element. So there is\n# no COM call sequence that can change that ordering, and none is needed:\n# the visible document content is unchanged.\n#\n# Reflect that faithfully as a no-op: read a harmless property and don't\n# mutate anything.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
